$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 109, pushing existing rows 109..171 down to 111..173.
$ws.Rows("109:110").Insert()

# --- New row 109 ---
$ws.Range("A109").Value = 4
$ws.Range("B109").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C109").Value = "Los Lagos"
$ws.Range("D109").Value = 44574
$ws.Range("E109").Value = 10
$ws.Range("F109").Value = 100112024
$ws.Range("G109").Value = "Choclo"
$ws.Range("H109").Value = "Choclero"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 4000
$ws.Range("K109").Value = 500
$ws.Range("L109").Value = 600
$ws.Range("M109").Value = 550
$ws.Range("N109").Value = "`$/unidad"
$ws.Range("O109").Value = "Región de O'Higgins"
$ws.Range("P109").Value = 550
$ws.Range("Q109").Value = 1
$ws.Range("R109").Value = "Hortaliza"

# --- New row 110 ---
$ws.Range("A110").Value = 4
$ws.Range("B110").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C110").Value = "Los Lagos"
$ws.Range("D110").Value = 44574
$ws.Range("E110").Value = 10
$ws.Range("F110").Value = 100112024
$ws.Range("G110").Value = "Choclo"
$ws.Range("H110").Value = "Choclero"
$ws.Range("I110").Value = "Segunda"
$ws.Range("J110").Value = 4000
$ws.Range("K110").Value = 350
$ws.Range("L110").Value = 350
$ws.Range("M110").Value = 350
$ws.Range("N110").Value = "`$/unidad"
$ws.Range("O110").Value = "Región de O'Higgins"
$ws.Range("P110").Value = 350
$ws.Range("Q110").Value = 1
$ws.Range("R110").Value = "Hortaliza"
